$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.873.46"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "1.758.15"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "327.87"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "0.4679"
$ws.Range("E7").Value = "  +1.26%  "
$ws.Range("D8").Value = "0.3504"
$ws.Range("D9").Value = "42.21"
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("D10").Value = "0.07353"
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("D11").Value = "1.080"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "20.49"
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").Value = "7.147"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").Value = "1.758.70"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "92.10"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").Value = "0.00001052"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").Value = "0.06405"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").Value = "1.0000"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "16.76"
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("D22").Value = "5.742"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").Value = "27.897.05"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("D25").Value = "2.152"
$ws.Range("E25").Value = "  +3.36%  "
$ws.Range("D26").Value = "162.13"
$ws.Range("E26").Value = "  -1.74%  "
$ws.Range("D27").Value = "19.98"
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("D28").Value = "1.958.06"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").Value = "2.164"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "122.60"
$ws.Range("E30").Value = "  -2.82%  "
$ws.Range("D31").Value = "1.069"
$ws.Range("E31").Value = "  -3.06%  "
$ws.Range("D32").Value = "0.09351"
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("D33").Value = "3.644"
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("D34").Value = "5.536"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").Value = "0.02262"
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("D37").Value = "0.06066"
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").Value = "0.2062"
$ws.Range("E38").Value = "  -1.41%  "
$ws.Range("D39").Value = "4.896"
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").Value = "0.6118"
$ws.Range("E40").Value = "  -2.91%  "
$ws.Range("D41").Value = "1.179"
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("D42").Value = "7.776"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").Value = "1.356"
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("D44").Value = "13.09"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("D45").Value = "3.731"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "0.5772"
$ws.Range("E46").Value = "  -1.99%  "
$ws.Range("D47").Value = "123.01"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").Value = "1.920"
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("D49").Value = "0.06803"
$ws.Range("E49").Value = "  -1.98%  "
$ws.Range("D50").Value = "1.121"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").Value = "71.98"
$ws.Range("E51").Value = "  -0.66%  "
